$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cola aqui os valores")

$ws.Range("B3").Value = "Língua Portuguesa"
$ws.Range("C3").Value = "Matemática"
$ws.Range("D3").Value = "Matemática"
$ws.Range("E3").Value = "História"
$ws.Range("I3").Value = "Matemática"
$ws.Range("J3").Value = "Inglês"
$ws.Range("K3").Value = "História"
$ws.Range("L3").Value = "Inglês"
$ws.Range("M3").Value = "Língua Portuguesa"
$ws.Range("B4").Value = "Inglês"
$ws.Range("C4").Value = "Matemática"
$ws.Range("D4").Value = "Matemática"
$ws.Range("E4").Value = "História"
$ws.Range("F4").Value = "Geografia"
$ws.Range("I4").Value = "Artes"
$ws.Range("J4").Value = "Música"
$ws.Range("K4").Value = "Inglês"
$ws.Range("L4").Value = "Língua Portuguesa"
$ws.Range("B5").Value = "Matemática"
$ws.Range("C5").Value = "Inglês"
$ws.Range("D5").Value = "Inglês"
$ws.Range("E5").Value = "Inglês"
$ws.Range("F5").Value = "Língua Portuguesa"
$ws.Range("I5").Value = "Artes"
$ws.Range("J5").Value = "Geografia"
$ws.Range("K5").Value = "Ed. Financeira"
$ws.Range("L5").Value = "Língua Portuguesa"
$ws.Range("M5").Value = "Inglês"
$ws.Range("C7").Value = "Ciências"
$ws.Range("D7").Value = "Língua Portuguesa"
$ws.Range("F7").Value = "Educação Física"
$ws.Range("J7").Value = "Língua Portuguesa"
$ws.Range("K7").Value = "Matemática"
$ws.Range("L7").Value = "Mind Makers"
$ws.Range("M7").Value = "Matemática"
$ws.Range("B8").Value = "Música"
$ws.Range("C8").Value = "Ciências"
$ws.Range("D8").Value = "Ensino Religioso"
$ws.Range("E8").Value = "Língua Portuguesa"
$ws.Range("F8").Value = "Mind Makers"
$ws.Range("J8").Value = "Língua Portuguesa"
$ws.Range("L8").Value = "Ensino Religioso"
$ws.Range("B13").Value = "Inglês"
$ws.Range("C13").Value = "Música"
$ws.Range("D13").Value = "Inglês"
$ws.Range("E13").Value = "Língua Portuguesa"
$ws.Range("F13").Value = "Ed. Financeira"
$ws.Range("I13").Value = "Geografia"
$ws.Range("J13").Value = "História"
$ws.Range("K13").Value = "Língua Portuguesa"
$ws.Range("L13").Value = "Matemática"
$ws.Range("M13").Value = "Inglês"
$ws.Range("B14").Value = "Língua Portuguesa"
$ws.Range("C14").Value = "Inglês"
$ws.Range("D14").Value = "Geografia"
$ws.Range("E14").Value = "Inglês"
$ws.Range("F14").Value = "Inglês"
$ws.Range("I14").Value = "Geografia"
$ws.Range("J14").Value = "História"
$ws.Range("K14").Value = "Língua Portuguesa"
$ws.Range("L14").Value = "Matemática"
$ws.Range("M14").Value = "Matemática"
$ws.Range("B15").Value = "Língua Portuguesa"
$ws.Range("C15").Value = "Língua Portuguesa"
$ws.Range("D15").Value = "Geografia"
$ws.Range("E15").Value = "História"
$ws.Range("I15").Value = "Inglês"
$ws.Range("J15").Value = "Música"
$ws.Range("K15").Value = "Língua Portuguesa"
$ws.Range("L15").Value = "Matemática"
$ws.Range("M15").Value = "Matemática"
$ws.Range("D17").Value = "Ensino Religioso"
$ws.Range("F17").Value = "Artes"
$ws.Range("J17").Value = "Artes"
$ws.Range("K17").Value = "Ciências"
$ws.Range("M17").Value = "Língua Portuguesa"
$ws.Range("C18").Value = "Matemática"
$ws.Range("D18").Value = "Língua Portuguesa"
$ws.Range("F18").Value = "Artes"
$ws.Range("J18").Value = "Artes"
$ws.Range("K18").Value = "Ciências"
$ws.Range("L18").Value = "Mind Makers"
$ws.Range("M18").Value = "Língua Portuguesa"
